# edit.ps1 - applies the two changes described by the target diff:
#
# 1. Slide 16's table (graphicFrame "Google Shape;213;p29", the 3rd shape on
#    the slide) gets its table style swapped from
#    {0A90250C-6C2E-452B-99DE-20A63DE1BAF1} to
#    {6CFBB9F0-18FA-4569-B892-F14FDEA3B197}.
#
# 2. The deck's main theme (ppt/theme/theme1.xml, used by the Slide Master)
#    swaps its 12-colour scheme from the "Integral" palette to the "Office"
#    palette (the colours theme2.xml - the Notes Master's theme - already
#    uses). The fontScheme/fmtScheme blocks of theme1.xml and theme2.xml are
#    already identical, so matching the colour scheme is the only
#    content-level change needed on the reachable theme part.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{6CFBB9F0-18FA-4569-B892-F14FDEA3B197}")

# --- 2. Theme colour scheme -------------------------------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# ThemeColorScheme.RGB uses the usual VBA RGB() packing (R + G*256 + B*65536).
# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colors.Item(1).RGB  = 0x000000    # dk1      000000 -> R=00 G=00 B=00
$colors.Item(2).RGB  = 0xFFFFFF    # lt1      FFFFFF -> R=FF G=FF B=FF
$colors.Item(3).RGB  = 0x6A5444    # dk2      44546A -> R=44 G=54 B=6A
$colors.Item(4).RGB  = 0xE6E6E7    # lt2      E7E6E6 -> R=E7 G=E6 B=E6
$colors.Item(5).RGB  = 0xD59B5B    # accent1  5B9BD5 -> R=5B G=9B B=D5
$colors.Item(6).RGB  = 0x317DED    # accent2  ED7D31 -> R=ED G=7D B=31
$colors.Item(7).RGB  = 0xA5A5A5    # accent3  A5A5A5 -> R=A5 G=A5 B=A5
$colors.Item(8).RGB  = 0x00C0FF    # accent4  FFC000 -> R=FF G=C0 B=00
$colors.Item(9).RGB  = 0xC47244    # accent5  4472C4 -> R=44 G=72 B=C4
$colors.Item(10).RGB = 0x47AD70    # accent6  70AD47 -> R=70 G=AD B=47
$colors.Item(11).RGB = 0xC16305    # hlink    0563C1 -> R=05 G=63 B=C1
$colors.Item(12).RGB = 0x724F95    # folHlink 954F72 -> R=95 G=4F B=72
